$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C32").Interior.ColorIndex = -4142

$styles = $wb.Styles
Write-Host "styles type:" $styles.GetType()
Write-Host "styles count:" $styles.Count
